{"js": "const replacements = [\n  [\"720\u00f76=120, 0\", \"147\u00f73=49, 0\"],\n  [\"627\u00f75=125, 2\", \"714\u00f72=357, 0\"],\n  [\"591\u00f74=147, 3\", \"265\u00f77=37, 6\"],\n  [\"897\u00f75=179, 2\", \"670\u00f77=95, 5\"],\n  [\"850\u00f75=170, 0\", \"598\u00f79=66, 4\"],\n  [\"720\u00f74=180, 0\", \"401\u00f75=80, 1\"],\n  [\"152\u00f79=16, 8\", \"592\u00f72=296, 0\"],\n  [\"606\u00f74=151, 2\", \"165\u00f75=33, 0\"],\n  [\"534\u00f79=59, 3\", \"206\u00f77=29, 3\"],\n  [\"658\u00f74=164, 2\", \"342\u00f78=42, 6\"],\n  [\"897\u00f79=99, 6\", \"734\u00f74=183, 2\"],\n  [\"460\u00f79=51, 1\", \"407\u00f78=50, 7\"],\n  [\"737\u00f74=184, 1\", \"273\u00f72=136, 1\"],\n  [\"617\u00f74=154, 1\", \"503\u00f79=55, 8\"],\n  [\"186\u00f72=93, 0\", \"809\u00f78=101, 1\"],\n  [\"566\u00f79=62, 8\", \"162\u00f79=18, 0\"],\n  [\"858\u00f77=122, 4\", \"868\u00f74=217, 0\"],\n  [\"465\u00f77=66, 3\", \"839\u00f77=119, 6\"],\n  [\"276\u00f75=55, 1\", \"267\u00f73=89, 0\"],\n  [\"279\u00f72=139, 1\", \"963\u00f75=192, 3\"],\n  [\"854\u00f77=122, 0\", \"418\u00f76=69, 4\"],\n  [\"719\u00f76=119, 5\", \"258\u00f77=36, 6\"],\n  [\"359\u00f75=71, 4\", \"178\u00f75=35, 3\"],\n  [\"863\u00f79=95, 8\", \"400\u00f78=50, 0\"],\n  [\"430\u00f76=71, 4\", \"411\u00f74=102, 3\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit division equation text with its updated value.\n$d = $word.ActiveDocument\n$replacements = @(\n    @{ Old = '720\u00f76=120, 0'; New = '147\u00f73=49, 0' },\n    @{ Old = '627\u00f75=125, 2'; New = '714\u00f72=357, 0' },\n    @{ Old = '591\u00f74=147, 3'; New = '265\u00f77=37, 6' },\n    @{ Old = '897\u00f75=179, 2'; New = '670\u00f77=95, 5' },\n    @{ Old = '850\u00f75=170, 0'; New = '598\u00f79=66, 4' },\n    @{ Old = '720\u00f74=180, 0'; New = '401\u00f75=80, 1' },\n    @{ Old = '152\u00f79=16, 8'; New = '592\u00f72=296, 0' },\n    @{ Old = '606\u00f74=151, 2'; New = '165\u00f75=33, 0' },\n    @{ Old = '534\u00f79=59, 3'; New = '206\u00f77=29, 3' },\n    @{ Old = '658\u00f74=164, 2'; New = '342\u00f78=42, 6' },\n    @{ Old = '897\u00f79=99, 6'; New = '734\u00f74=183, 2' },\n    @{ Old = '460\u00f79=51, 1'; New = '407\u00f78=50, 7' },\n    @{ Old = '737\u00f74=184, 1'; New = '273\u00f72=136, 1' },\n    @{ Old = '617\u00f74=154, 1'; New = '503\u00f79=55, 8' },\n    @{ Old = '186\u00f72=93, 0'; New = '809\u00f78=101, 1' },\n    @{ Old = '566\u00f79=62, 8'; New = '162\u00f79=18, 0' },\n    @{ Old = '858\u00f77=122, 4'; New = '868\u00f74=217, 0' },\n    @{ Old = '465\u00f77=66, 3'; New = '839\u00f77=119, 6' },\n    @{ Old = '276\u00f75=55, 1'; New = '267\u00f73=89, 0' },\n    @{ Old = '279\u00f72=139, 1'; New = '963\u00f75=192, 3' },\n    @{ Old = '854\u00f77=122, 0'; New = '418\u00f76=69, 4' },\n    @{ Old = '719\u00f76=119, 5'; New = '258\u00f77=36, 6' },\n    @{ Old = '359\u00f75=71, 4'; New = '178\u00f75=35, 3' },\n    @{ Old = '863\u00f79=95, 8'; New = '400\u00f78=50, 0' },\n    @{ Old = '430\u00f76=71, 4'; New = '411\u00f74=102, 3' }\n)\n\nforeach ($rep in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $rep.Old\n    $find.Replacement.Text = $rep.New\n    $find.Execute($rep.Old, $false, $false, $false, $false, $false, $true, 1, $false, $rep.New, 2) | Out-Null\n}\n"}
